# Suite.xlsx edit: drop BankManagerSuite/CustomerSuite rows, keep only the
# CarWaleSuite row (row 2), clearing rows 3 and 4 back to blank cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 already reads CarWaleSuite / Y in most saves of this sheet, but set
# it explicitly so the result is correct regardless of starting state.
$ws.Range("A2").Value = "CarWaleSuite"
$ws.Range("B2").Value = "Y"

# Rows 3-4 (BankManagerSuite/Y and CustomerSuite/Y) are removed, leaving the
# cells blank (formatting/style stays put).
$ws.Range("A3:B4").ClearContents()

# Matches the saved selection in the edited workbook.
$ws.Range("A3:B4").Select()
